# "triying out something new"
#
# Re-creates a short PowerPoint edit session:
#   1. Nudge/resize the title placeholder on slide 1 (text stays "hahaha").
#   2. Add a new second slide using the "Title and Content" layout and
#      type a title on it ("haaaaaaaaaa"), leaving the content placeholder
#      empty.

$p = $ppt.ActivePresentation

# --- 1. Slide 1: reposition/resize the ctrTitle placeholder -------------
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1)
$title1.Left = 119.225827
$title1.Top = 88.375045
$title1.Width = 720.0
$title1.Height = 188.0

# --- 2. Insert a new slide after slide 1 ---------------------------------
$ppLayoutText = 2
$newSlide = $p.Slides.Add($s1.SlideIndex + 1, $ppLayoutText)

$newTitle = $newSlide.Shapes.Item(1)
$newTitle.TextFrame.TextRange.Text = "haaaaaaaaaa"
